# feat: added vba code to format pdf and added password to pdf
# Re-derive the BRL unit-cost / total-price columns (I, J) from the raw
# "<amount> <CCY>" strings in columns G/H (strip the currency code), and
# bump the exchange-conversion date in column K to the new run date.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count
if (-not $lastRow -or $lastRow -lt 2) { $lastRow = 30 }

for ($r = 2; $r -le $lastRow; $r++) {

    $gVal = $ws.Cells.Item($r, 7).Value2   # column G - "UNIT COST"
    $hVal = $ws.Cells.Item($r, 8).Value2   # column H - "TOTAL PRICE"
    $iOld = $ws.Cells.Item($r, 9).Value2   # column I - "Unit Cost (BRL)" (pre-edit)
    $jOld = $ws.Cells.Item($r, 10).Value2  # column J - "Total Price (BRL)" (pre-edit)

    if ($gVal -and $iOld) {
        $amount = ($gVal -split ' ')[0]
        $cell = $ws.Cells.Item($r, 9)      # column I - "Unit Cost (BRL)"
        $cell.NumberFormat = "@"
        $cell.Value = $amount
    }

    if ($hVal -and $jOld) {
        $amount = ($hVal -split ' ')[0]
        $cell = $ws.Cells.Item($r, 10)     # column J - "Total Price (BRL)"
        $cell.NumberFormat = "@"
        $cell.Value = $amount
    }

    $kCell = $ws.Cells.Item($r, 11)        # column K - "Exchange Conversion Date/Time"
    $kCell.NumberFormat = "@"
    $kCell.Value = "05/08/2025"
}
